$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "struggle", 2.563363254070282, -5.564052700996399, -4.925167679786682, -0.0056505035609006, -0.007177666760981, 0),
    @(100, "struggle", 2.571200489997864, -5.45090651512146, -4.94497549533844, -0.0161879286170005, 0.0122173046693205, -0.0047342055477201),
    @(200, "struggle", 2.582025349140167, -5.429405391216278, -4.891633093357086, 0.0029016099870204, -0.0010690141934901, -0.0094684110954403),
    @(300, "struggle", 2.521161556243896, -5.436496257781982, -4.74793529510498, 0.00534507073462, 0.0088575463742017, 0.0045814891345798),
    @(400, "struggle", 2.451793313026428, -5.718138635158539, -4.555151760578156, 0.0036651915870606, -0.0503963828086853, 0.0389426611363887),
    @(500, "struggle", 4.312116026878357, -5.115874052047729, -7.913362622261047, -1.588860511779785, -2.905885934829712, 0.6783658266067505),
    @(600, "struggle", 3.11747863888741, -5.436717748641967, -5.15502381324769, -0.9424123764038086, -3.084411382675171, 0.4729624092578888),
    @(700, "struggle", 3.307996869087219, -5.476139068603516, -0.757482767105103, 0.2102903574705124, -1.237612962722778, 0.2101376503705978),
    @(800, "struggle", 4.280841529369356, -4.743481069803241, 2.116865754127498, -0.4138612151145935, -0.9240863919258118, 0.252134621143341),
    @(900, "struggle", 1.148419260978671, -4.571187555789954, 3.019693136215206, 0.432645320892334, 1.750434398651123, -0.1994474977254867),
    @(1000, "struggle", -5.23933637142183, -7.058503150939933, -0.1314393877982891, 1.822669148445129, 5.886908531188965, -2.078774452209473),
    @(1100, "struggle", 5.029654502868652, -12.30031204223633, -12.83510589599609, 0.7883216142654419, 2.360841512680054, 1.160491228103638),
    @(1200, "struggle", -1.211387872695933, -5.499476730823506, -6.533813059329978, 0.5829181671142578, 1.261283993721008, 0.6953173875808716),
    @(1300, "struggle", -3.101414084434511, -3.493300497531889, -7.406978726387027, -0.3068070709705353, 1.792431354522705, -0.2423607856035232),
    @(1400, "struggle", -3.497859179973603, -3.212260961532593, -10.54216539859774, 0.2667953968048095, -0.7637342810630798, 0.0548251569271087),
    @(1500, "struggle", -1.64927089214325, -7.43125307559967, -8.658325910568239, -1.207069754600525, -6.127589225769043, -1.29213273525238),
    @(1600, "struggle", -1.025534451007847, -11.15312564373016, -2.677394092082968, -1.138041973114014, -3.366478443145752, 0.9393580555915833),
    @(1700, "struggle", -1.654657959938065, -9.341210365295419, 0.4500467777252322, 1.038929104804993, -0.7336491346359253, -1.004567861557007),
    @(1800, "struggle", 2.828480809926987, -5.927687406539919, -2.06907200813294, -0.5413793325424194, -1.692554831504822, -0.2292271852493286),
    @(1900, "struggle", 3.216025352478028, -4.024631500244141, 4.593602180480957, -0.3520110845565796, 0.6711881756782532, -0.102472648024559),
    @(2000, "struggle", -1.550865292549141, -5.700099587440493, 3.356234908103941, 2.986520290374756, 4.694652080535889, 0.0166460778564214),
    @(2100, "struggle", -2.655673027038576, -3.905611395835868, -3.318085908889795, 1.599092483520508, 4.835456371307373, 1.327257513999939),
    @(2200, "struggle", -1.430848956108096, -1.578429281711582, -8.26511919498443, -0.3246748745441437, 0.3605632185935974, 0.3962988257408142),
    @(2300, "struggle", -2.407764434814463, -1.766093373298651, -8.389460563659659, 0.05283984541893, 0.2293798923492431, -0.2738203406333923),
    @(2400, "struggle", -2.611050009727472, -2.573673054575919, -8.294337868690496, -0.6475171446800232, -0.2813034355640411, 0.0232128798961639),
    @(2500, "struggle", 0.61596310138703, -2.870795279741297, -8.901223957538587, -0.1372919678688049, -2.705674886703491, -0.5198463201522827),
    @(2600, "struggle", 2.079445004463198, -5.368536770343783, -3.524431616067886, -0.6478226184844971, -0.2128865420818328, -0.0656680166721344),
    @(2700, "struggle", 0.3216586112976074, -3.676267147064209, -3.865855693817138, -0.1000291854143142, 0.1372919678688049, -0.1838704347610473),
    @(2800, "struggle", 1.706132471561434, -4.47040206193924, -5.197765350341799, 0.2654209434986114, 0.0520762614905834, 0.0438295826315879),
    @(2900, "struggle", 1.816254138946533, -3.920204520225524, -5.85852086544037, 0.0809396430850029, 0.3397938013076782, 0.0826195254921913)
)

$startRow = 2
$numRows = $data.Count
$numCols = 8

for ($i = 0; $i -lt $numRows; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    for ($j = 0; $j -lt $numCols; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rec[$j]
    }
}
